$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2
$ws.Range("A2").Value = "suka"
$ws.Range("B2").Value = "blyad'"
$ws.Range("C2").Value = 32

# Add new row 3
$ws.Range("A3").Value = "suka"
$ws.Range("B3").Value = "pidr"
$ws.Range("C3").Value = 232
